$wb = $excel.ActiveWorkbook

# This script applies refreshed market-price/profit figures to several
# rows across multiple sheets, as pulled by the scheduled data runner.

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 2811.1562
$ws.Range("I64").Value = 2778.5715
$ws.Range("J64").Value = 3039.25
$ws.Range("K64").Value = 2778.5715
$ws.Range("L64").Value = 3039.25
$ws.Range("M64").Value = -2530.5715
$ws.Range("N64").Value = -3535.25

$ws.Range("H67").Value = 2811.1562
$ws.Range("I67").Value = 2778.5715
$ws.Range("J67").Value = 3039.25
$ws.Range("K67").Value = 2778.5715
$ws.Range("L67").Value = 3039.25
$ws.Range("M67").Value = -1920.5715
$ws.Range("N67").Value = -4755.25

$ws.Range("H74").Value = 9793.333000000001
$ws.Range("I74").Value = 12590
$ws.Range("J74").Value = 4200
$ws.Range("K74").Value = 12590
$ws.Range("L74").Value = 4200
$ws.Range("M74").Value = -11654
$ws.Range("N74").Value = -6072

$ws.Range("H77").Value = 9793.333000000001
$ws.Range("I77").Value = 12590
$ws.Range("J77").Value = 4200
$ws.Range("K77").Value = 62950
$ws.Range("L77").Value = 21000
$ws.Range("M77").Value = -58270
$ws.Range("N77").Value = -30360

$ws.Range("H98").Value = 1908.75
$ws.Range("I98").Value = 1190.5
$ws.Range("K98").Value = 1190.5
$ws.Range("M98").Value = 307.5

$ws.Range("H106").Value = 2182.9167
$ws.Range("I106").Value = 2017.7273
$ws.Range("J106").Value = 4000
$ws.Range("K106").Value = 2017.7273
$ws.Range("L106").Value = 4000
$ws.Range("M106").Value = -1386.7273
$ws.Range("N106").Value = -5262

$ws.Range("H122").Value = 1908.75
$ws.Range("I122").Value = 1190.5
$ws.Range("K122").Value = 3571.5
$ws.Range("M122").Value = -1121.5

$ws.Range("H135").Value = 19978.623
$ws.Range("I135").Value = 26080.574
$ws.Range("J135").Value = 1203.3846
$ws.Range("K135").Value = 234725.166
$ws.Range("L135").Value = 10830.4614
$ws.Range("M135").Value = -232190.166
$ws.Range("N135").Value = -15900.4614

$ws.Range("H137").Value = 3489574.2
$ws.Range("I137").Value = 1390035.8
$ws.Range("J137").Value = 14287200
$ws.Range("K137").Value = 4170107.4
$ws.Range("L137").Value = 42861600
$ws.Range("M137").Value = -4167557.4
$ws.Range("N137").Value = -42866700

$ws.Range("H138").Value = 2404.3137
$ws.Range("I138").Value = 2761.9167
$ws.Range("J138").Value = 2294.282
$ws.Range("K138").Value = 8285.750100000001
$ws.Range("L138").Value = 6882.846
$ws.Range("M138").Value = -3145.750100000001
$ws.Range("N138").Value = -17162.846

$ws.Range("H141").Value = 1819.8864
$ws.Range("I141").Value = 1435.5
$ws.Range("J141").Value = 2643.5715
$ws.Range("K141").Value = 4306.5
$ws.Range("L141").Value = 7930.7145
$ws.Range("M141").Value = 873.5
$ws.Range("N141").Value = -18290.7145

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H55").Value = 16049.667
$ws.Range("J55").Value = 16049.667
$ws.Range("L55").Value = 16049.667
$ws.Range("N55").Value = -16679.667

$ws.Range("H74").Value = 4348704
$ws.Range("I74").Value = 868.9286
$ws.Range("K74").Value = 868.9286
$ws.Range("M74").Value = 5.07140000000004

$ws.Range("H77").Value = 4348704
$ws.Range("I77").Value = 868.9286
$ws.Range("K77").Value = 4344.643
$ws.Range("M77").Value = 23.35699999999997

$ws.Range("H110").Value = 10854.45
$ws.Range("I110").Value = 10886.823
$ws.Range("J110").Value = 10671
$ws.Range("K110").Value = 10886.823
$ws.Range("L110").Value = 10671
$ws.Range("M110").Value = -8841.823
$ws.Range("N110").Value = -14761

$ws.Range("H132").Value = 141433.89
$ws.Range("I132").Value = 157672.56
$ws.Range("K132").Value = 473017.68
$ws.Range("M132").Value = -470487.68

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 53000.61
$ws.Range("J82").Value = 30710.525
$ws.Range("L82").Value = 30710.525
$ws.Range("N82").Value = -31476.525

$ws.Range("H85").Value = 53000.61
$ws.Range("J85").Value = 30710.525
$ws.Range("L85").Value = 30710.525
$ws.Range("N85").Value = -33362.525

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1958.4375
$ws.Range("I31").Value = 1471.6522
$ws.Range("J31").Value = 3202.4443
$ws.Range("K31").Value = 1471.6522
$ws.Range("L31").Value = 3202.4443
$ws.Range("M31").Value = -1176.6522
$ws.Range("N31").Value = -3792.4443

$ws.Range("H34").Value = 1958.4375
$ws.Range("I34").Value = 1471.6522
$ws.Range("J34").Value = 3202.4443
$ws.Range("K34").Value = 1471.6522
$ws.Range("L34").Value = 3202.4443
$ws.Range("M34").Value = -1269.6522
$ws.Range("N34").Value = -3606.4443

$ws.Range("H58").Value = 864.04443
$ws.Range("I58").Value = 860.7692
$ws.Range("K58").Value = 860.7692
$ws.Range("M58").Value = -657.7692

$ws.Range("H132").Value = 1713.6383
$ws.Range("I132").Value = 1611.575
$ws.Range("K132").Value = 4834.725
$ws.Range("M132").Value = -2304.725

$ws.Range("H136").Value = 864.04443
$ws.Range("I136").Value = 860.7692
$ws.Range("K136").Value = 2582.3076
$ws.Range("M136").Value = -32.30760000000009

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 559.73334
$ws.Range("I14").Value = 559.73334
$ws.Range("K14").Value = 1679.20002
$ws.Range("M14").Value = -1506.20002

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 3169.0476
$ws.Range("I122").Value = 3686
$ws.Range("J122").Value = 2699.0908
$ws.Range("K122").Value = 11058
$ws.Range("L122").Value = 8097.2724
$ws.Range("M122").Value = -8608
$ws.Range("N122").Value = -12997.2724

$ws.Range("H132").Value = 2144
$ws.Range("I132").Value = 1704.6957
$ws.Range("J132").Value = 2705.3333
$ws.Range("K132").Value = 5114.0871
$ws.Range("L132").Value = 8115.999899999999
$ws.Range("M132").Value = -2584.0871
$ws.Range("N132").Value = -13175.9999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 2556.639
$ws.Range("I132").Value = 2072.7144
$ws.Range("J132").Value = 4250.375
$ws.Range("K132").Value = 6218.1432
$ws.Range("L132").Value = 12751.125
$ws.Range("M132").Value = -3688.1432
$ws.Range("N132").Value = -17811.125

$ws.Range("H133").Value = 42705.43
$ws.Range("J133").Value = 42705.43
$ws.Range("L133").Value = 42705.43
$ws.Range("N133").Value = -47765.43

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H109").Value = 16917.7
$ws.Range("J109").Value = 16917.7
$ws.Range("L109").Value = 16917.7
$ws.Range("N109").Value = -19691.7

$ws.Range("H132").Value = 1172.1846
$ws.Range("I132").Value = 1097.3103
$ws.Range("J132").Value = 1792.5714
$ws.Range("K132").Value = 3291.9309
$ws.Range("L132").Value = 5377.7142
$ws.Range("M132").Value = -761.9309000000003
$ws.Range("N132").Value = -10437.7142
